# Adds the new "AAA" (Administrator users) common-code rows used by the
# notifications migration, and a couple of related CCP limit-checker fixes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended right after the existing data (rows 1-173 already in use).
# Columns: A=COMM1_CD  B=COMM2_CD  C=COMM_NM  D=REG_ID  E=REG_DTM

$ws.Cells.Item(174, 1).Value = "AAA"
$ws.Cells.Item(174, 2).Value = "$$"
$ws.Cells.Item(174, 3).Value = "Administrator users"
$ws.Cells.Item(174, 4).Value = "admin"
$ws.Cells.Item(174, 5).Value = 20201201153327

$ws.Cells.Item(175, 1).Value = "AAA"
$ws.Cells.Item(175, 2).Value = "admin"
$ws.Cells.Item(175, 3).Value = "Default Administrator"
$ws.Cells.Item(175, 4).Value = "admin"
$ws.Cells.Item(175, 5).Value = 20201201153327

$ws.Cells.Item(176, 1).Value = "AAA"
$ws.Cells.Item(176, 2).Value = "kenny"
$ws.Cells.Item(176, 3).Value = "Kwon Yoon"
$ws.Cells.Item(176, 4).Value = "admin"
$ws.Cells.Item(176, 5).Value = 20201201153327

# Match the formatting used by the preceding block of rows (163-173): column A
# uses style index 4 (the Dotum/Korean font used for the "AAA" group header
# column), columns D/E use style index 1 (already the default via copy below).
$ws.Range("A163").Copy()
$ws.Range("A174:A176").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D163:E163").Copy()
$ws.Range("D174:E176").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A174:E176").RowHeight = 15.75

$excel.CutCopyMode = 0
